$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_registros values for rows that keep the same empadronador name
$ws.Range("B2").Value = 171
$ws.Range("B3").Value = 160
$ws.Range("B4").Value = 134
$ws.Range("B5").Value = 134
$ws.Range("B6").Value = 134
$ws.Range("B8").Value = 121
$ws.Range("B9").Value = 118

# Rows 7, 10, 11 have their empadronador names re-ordered along with new totals
$ws.Range("A7").Value = "CHIROQUE YARLEQUE BETTY ELIZABETH"
$ws.Range("B7").Value = 123

$ws.Range("A10").Value = "HERRERA JUAN MANUEL"
$ws.Range("B10").Value = 116

$ws.Range("A11").Value = "SILVA ALVARADO EVELYN DE JESUS"
$ws.Range("B11").Value = 110
